$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data: Letter "D", Number 4
$ws.Range("A5").Value = "D"
$ws.Range("B5").Value = 4

# Update the selection to match the post-edit state recorded in the workbook
$ws.Range("D8").Select()
